{"js": "// The template's header carries a labelled placeholder line:\n//   Date for interview: \"DD/MM/YYYY\"\n// The field is actually used to record the date the participant signs the\n// consent form, so \"interview\" is corrected to \"signature\":\n//   Date for signature: \"DD/MM/YYYY\"\n// The label lives in the document header, not the body, so every section's\n// primary header is searched for the phrase and updated in place.\n\nconst sections = context.document.sections;\nsections.load(\"items\");\nawait context.sync();\n\nfor (const section of sections.items) {\n  const header = section.getHeader(Word.HeaderFooterType.primary);\n\n  const matches = header.search(\"Date for interview\", { matchCase: true });\n  matches.load(\"items/text\");\n  await context.sync();\n\n  for (const match of matches.items) {\n    match.insertText(\"Date for signature\", Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# The template's header carries a labelled placeholder line:\n#   Date for interview: \"DD/MM/YYYY\"\n# The field actually records the date the participant signs the consent\n# form, so \"interview\" is corrected to \"signature\":\n#   Date for signature: \"DD/MM/YYYY\"\n# The label lives in the document header (not the body), so every\n# section's primary header range is searched for the phrase and updated\n# in place via Find/Replace, which preserves the surrounding run\n# formatting (font, color, size, language) untouched.\n\n$d = $word.ActiveDocument\n\nforeach ($sec in $d.Sections) {\n    $hdr = $sec.Headers.Item([Microsoft.Office.Interop.Word.WdHeaderFooterIndex]::wdHeaderFooterPrimary)\n    if ($hdr.Exists) {\n        $rng = $hdr.Range\n        $find = $rng.Find\n        $find.ClearFormatting()\n        $find.Replacement.ClearFormatting()\n        $find.Execute(\n            \"Date for interview\",  # FindText\n            $true,                 # MatchCase\n            $false,                # MatchWholeWord\n            $false,                # MatchWildcards\n            $false,                # MatchSoundsLike\n            $false,                # MatchAllWordForms\n            $true,                 # Forward\n            1,                     # Wrap: wdFindContinue\n            $false,                # Format\n            \"Date for signature\",  # ReplaceWith\n            2                      # Replace: wdReplaceAll\n        )\n    }\n}\n"}
